$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the APS (Annual Population Survey) period data - the "latest period" and
# "next period" values were stale / pointing at a period with no data yet.
# Shift Latest period (col C) and Next period (col D) forward for the three
# "Employment ..." rows (2-4), and fix the qualifications row (7) so its
# "Next period" points at the correct upcoming period instead of a period
# that has no data.

$newLatest = "Apr 2023 - Mar 2024 (18/07/24)"
$newNextCommon = "Jul 2023 - Jun 2024 (15/10/24)"
$newNextQualifications = "Jan 2024 - Dec 2024 (TBC)"

# Rows 2-4: Employment volumes, Employment by occupation, Employment by industry
$ws.Range("C2").Value = $newLatest
$ws.Range("D2").Value = $newNextCommon

$ws.Range("C3").Value = $newLatest
$ws.Range("D3").Value = $newNextCommon

$ws.Range("C4").Value = $newLatest
$ws.Range("D4").Value = $newNextCommon

# Row 7: Highest qualification level by age and gender - latest period (C7)
# stays the same, only the next period (D7) needs correcting since the
# previously listed period had no available data.
$ws.Range("D7").Value = $newNextQualifications

# Restore the user's active cell selection to D2, as left after the edit.
$ws.Range("D2").Select()
